# Apply the edits described by the diff:
# 1. Update the letter date from September 19, 2025 to September 21, 2025.
# 2. Split the single "1754 Los Altos Drive, San Mateo CA 94402" address line
#    into two separate paragraphs ("1754 Los Altos Drive" and
#    "San Mateo, CA 94402") and add a new blank paragraph after it.
# 3. Remove the two blank paragraphs that followed "...Board of Directors".

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# --- 1. Update the date -----------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# --- 2. Split the mailing address paragraph ----------------------------
$addrIndex = Find-ParagraphIndex $d "1754 Los Altos Drive"
$addrPara = $d.Paragraphs.Item($addrIndex)

# Replace the paragraph's text (everything except its trailing paragraph
# mark) with the street line, a paragraph break, the city/state/zip line
# and another paragraph break. This produces: "1754 Los Altos Drive" /
# "San Mateo, CA 94402" / <new blank paragraph>, immediately followed by
# the blank paragraph that was already there.
$addrBody = $d.Range($addrPara.Range.Start, $addrPara.Range.End - 1)
$addrBody.Text = "1754 Los Altos Drive" + [char]13 + "San Mateo, CA 94402" + [char]13

# --- 3. Remove the two blank paragraphs after "Board of Directors" -----
$bodIndex = Find-ParagraphIndex $d "Board of Directors"

# Deleting the (empty) paragraph's range removes just its paragraph mark,
# merging it with what follows. Doing this twice removes both blank
# paragraphs that sat between "Board of Directors" and the next real
# paragraph.
$d.Paragraphs.Item($bodIndex + 1).Range.Delete()
$d.Paragraphs.Item($bodIndex + 1).Range.Delete()
